$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5714285714285714
$ws.Range("J2").Value = 0.02164502164502164
$ws.Range("P2").Value = 0.1341991341991342
$ws.Range("S2").Value = 0.08225108225108226
$ws.Range("C3").Value = 0.02205882352941177
$ws.Range("P3").Value = 0.7279411764705882
$ws.Range("S3").Value = 0.25
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.6486486486486487
$ws.Range("S4").Value = 0.2702702702702703
$ws.Range("B6").Value = 0.05581395348837209
$ws.Range("D6").Value = 0.004651162790697674
$ws.Range("F6").Value = 0.05581395348837209
$ws.Range("J6").Value = 0.2604651162790698
$ws.Range("O6").Value = 0.01395348837209302
$ws.Range("Q6").Value = 0.172093023255814
$ws.Range("R6").Value = 0.05116279069767442
$ws.Range("S6").Value = 0.386046511627907
$ws.Range("B7").Value = 0.04705882352941176
$ws.Range("F7").Value = 0.07058823529411765
$ws.Range("J7").Value = 0.1529411764705882
$ws.Range("O7").Value = 0.01764705882352941
$ws.Range("Q7").Value = 0.1705882352941177
$ws.Range("R7").Value = 0.08235294117647059
$ws.Range("S7").Value = 0.4588235294117647
$ws.Range("B8").Value = 0.06425702811244979
$ws.Range("D8").Value = 0.01807228915662651
$ws.Range("F8").Value = 0.07228915662650602
$ws.Range("J8").Value = 0.1184738955823293
$ws.Range("O8").Value = 0.01606425702811245
$ws.Range("Q8").Value = 0.1907630522088354
$ws.Range("R8").Value = 0.07429718875502007
$ws.Range("S8").Value = 0.4457831325301205
$ws.Range("B9").Value = 0.04294478527607362
$ws.Range("D9").Value = 0.01226993865030675
$ws.Range("F9").Value = 0.049079754601227
$ws.Range("J9").Value = 0.147239263803681
$ws.Range("O9").Value = 0.006134969325153374
$ws.Range("Q9").Value = 0.1779141104294479
$ws.Range("R9").Value = 0.049079754601227
$ws.Range("S9").Value = 0.5153374233128835
$ws.Range("B10").Value = 0.09912767644726407
$ws.Range("D10").Value = 0.02141157811260904
$ws.Range("F10").Value = 0.06661379857256146
$ws.Range("J10").Value = 0.119746233148295
$ws.Range("O10").Value = 0.006344171292624901
$ws.Range("Q10").Value = 0.2307692307692308
$ws.Range("R10").Value = 0.06582077716098335
$ws.Range("S10").Value = 0.3901665344964314
$ws.Range("G11").Value = 0.1371428571428571
$ws.Range("J11").Value = 0.1457142857142857
$ws.Range("K11").Value = 0.24
$ws.Range("L11").Value = 0.4514285714285714
$ws.Range("S11").Value = 0.02571428571428571
$ws.Range("G12").Value = 0.6645962732919255
$ws.Range("J12").Value = 0.2608695652173913
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.02484472049689441
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.5882352941176471
$ws.Range("J13").Value = 0.3529411764705883
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.03061224489795918
$ws.Range("H15").Value = 0.2040816326530612
$ws.Range("I15").Value = 0.05102040816326531
$ws.Range("J15").Value = 0.2959183673469388
$ws.Range("K15").Value = 0.05612244897959184
$ws.Range("M15").Value = 0.01020408163265306
$ws.Range("O15").Value = 0.07653061224489796
$ws.Range("S15").Value = 0.2755102040816326
$ws.Range("F16").Value = 0.02649006622516556
$ws.Range("H16").Value = 0.1721854304635762
$ws.Range("I16").Value = 0.04635761589403974
$ws.Range("J16").Value = 0.4172185430463576
$ws.Range("K16").Value = 0.08609271523178808
$ws.Range("M16").Value = 0.01986754966887417
$ws.Range("O16").Value = 0.04635761589403974
$ws.Range("S16").Value = 0.1854304635761589
$ws.Range("F17").Value = 0.01455301455301455
$ws.Range("H17").Value = 0.1912681912681913
$ws.Range("I17").Value = 0.09355509355509356
$ws.Range("J17").Value = 0.3887733887733888
$ws.Range("K17").Value = 0.09563409563409564
$ws.Range("M17").Value = 0.01247401247401247
$ws.Range("N17").Value = 0.002079002079002079
$ws.Range("O17").Value = 0.04573804573804574
$ws.Range("S17").Value = 0.1559251559251559
$ws.Range("F18").Value = 0.02666666666666667
$ws.Range("H18").Value = 0.16
$ws.Range("I18").Value = 0.1066666666666667
$ws.Range("J18").Value = 0.38
$ws.Range("K18").Value = 0.08666666666666667
$ws.Range("M18").Value = 0.02
$ws.Range("O18").Value = 0.04
$ws.Range("S18").Value = 0.18
$ws.Range("F19").Value = 0.01468531468531468
$ws.Range("H19").Value = 0.2237762237762238
$ws.Range("I19").Value = 0.06013986013986014
$ws.Range("J19").Value = 0.3412587412587413
$ws.Range("K19").Value = 0.1237762237762238
$ws.Range("M19").Value = 0.01468531468531468
$ws.Range("N19").Value = 0.0006993006993006993
$ws.Range("O19").Value = 0.06573426573426573
$ws.Range("S19").Value = 0.1552447552447553
